# Apply the changes described in the commit:
# "Add text to historical and update baseline sources."
#
# 1) Populate column D (source note) for rows 57-70 with the appropriate
#    "source" shared strings (re-using existing strings where possible,
#    and introducing one brand-new string for row 66).
# 2) Bump row 66's height to fit the newly-added (longer) text.
# 3) Refresh workbook-level bookkeeping (author path in absPath,
#    revision id) and view state (window size, zoom, scroll position,
#    selected cell) to reflect the editor's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column D additions for rows 57-70 -------------------------------
$ws.Range("D57").Value = "Developed by RAPID Team"
$ws.Range("D58").Value = "Developed by RAPID Team"
$ws.Range("D59").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D60").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D61").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D62").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D63").Value = "RAPID Team Modified from U.S. Census "
$ws.Range("D64").Value = "Developed by RAPID Team"
$ws.Range("D65").Value = "Developed by RAPID Team"
$ws.Range("D66").Value = "RAPID Team Modified from National Compensation Survey"
$ws.Range("D67").Value = "Developed by RAPID Team"
$ws.Range("D68").Value = "RAPID Team Modified"
$ws.Range("D69").Value = "Developed by RAPID Team"
$ws.Range("D70").Value = "Developed by RAPID Team"

# --- 2) Row 66 grew taller to fit its new text --------------------------
$ws.Rows.Item(66).RowHeight = 96

# --- 3) Workbook / view bookkeeping -------------------------------------
$wb.ActiveSheet.Activate()
$ws.Range("D69").Select()
$excel.ActiveWindow.Zoom = 125
